# Updates cryptos list price (D) and 1h volume change (E) columns for rows 2-51.
# Values are plain text (not numbers/percentages) in the source data, so the
# helper below forces text formatting before the write and then clears the
# formatting again (so no stray number-format style is left on the cell) --
# this avoids Excel auto-converting numeric-looking strings (e.g. "355.67")
# into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellref, $val)
    $r = $ws.Range($cellref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue $ws 'D2' '51.823.09'
$ws.Range('E2').Value = '  -0.49%  '
Set-TextValue $ws 'D3' '2.928.39'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue $ws 'D5' '355.67'
$ws.Range('E5').Value = '  +1.09%  '
Set-TextValue $ws 'D6' '110.91'
$ws.Range('E6').Value = '  -0.54%  '
Set-TextValue $ws 'D7' '0.567'
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('E8').Value = '  +0.03%  '
Set-TextValue $ws 'D9' '0.628'
$ws.Range('E9').Value = '  +0.80%  '
Set-TextValue $ws 'D10' '39.37'
$ws.Range('E10').Value = '  -1.21%  '
Set-TextValue $ws 'D11' '0.0881'
$ws.Range('E11').Value = '  +2.92%  '
Set-TextValue $ws 'D12' '0.137'
$ws.Range('E12').Value = '  +0.63%  '
Set-TextValue $ws 'D13' '19.73'
$ws.Range('E13').Value = '  -1.11%  '
Set-TextValue $ws 'D14' '7.91'
$ws.Range('E14').Value = '  +1.83%  '
Set-TextValue $ws 'D15' '3.389.46'
$ws.Range('E15').Value = '  +1.25%  '
Set-TextValue $ws 'D16' '2.931.69'
$ws.Range('E16').Value = '  +0.87%  '
Set-TextValue $ws 'D17' '0.986'
$ws.Range('E17').Value = '  -2.12%  '
Set-TextValue $ws 'D18' '51.860.00'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('E19').Value = '  -1.57%  '
Set-TextValue $ws 'D20' '7.56'
$ws.Range('E20').Value = '  -2.03%  '
Set-TextValue $ws 'D21' '14.05'
$ws.Range('E21').Value = '  -1.88%  '
Set-TextValue $ws 'D22' '0.0₃0983'
$ws.Range('E22').Value = '  +0.26%  '
Set-TextValue $ws 'D23' '70.95'
$ws.Range('E23').Value = '  +0.33%  '
Set-TextValue $ws 'D24' '270.63'
Set-TextValue $ws 'D25' '2.82'
$ws.Range('E25').Value = '  +1.34%  '
Set-TextValue $ws 'D26' '0.184'
$ws.Range('E26').Value = '  +12.63%  '
Set-TextValue $ws 'D27' '27.19'
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('E28').Value = '  +0.21%  '
Set-TextValue $ws 'D29' '7.42'
$ws.Range('E29').Value = '  +15.38%  '
$ws.Range('E30').Value = '  +12.57%  '
Set-TextValue $ws 'D31' '39.20'
$ws.Range('E31').Value = '  +2.05%  '
Set-TextValue $ws 'D32' '10.59'
$ws.Range('E32').Value = '  +0.71%  '
Set-TextValue $ws 'D33' '6.07'
$ws.Range('E33').Value = '  -1.33%  '
Set-TextValue $ws 'D34' '52.21'
$ws.Range('E34').Value = '  -1.50%  '
Set-TextValue $ws 'D35' '0.0444'
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -14.89%  '
Set-TextValue $ws 'D38' '3.25'
$ws.Range('E38').Value = '  -1.82%  '
Set-TextValue $ws 'D39' '18.56'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  -1.35%  '
Set-TextValue $ws 'D41' '2.76'
$ws.Range('E41').Value = '  +4.06%  '
$ws.Range('E42').Value = '  +2.74%  '
$ws.Range('E43').Value = '  +2.18%  '
Set-TextValue $ws 'D44' '119.40'
$ws.Range('E44').Value = '  -2.17%  '
Set-TextValue $ws 'D45' '2.16'
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('E47').Value = '  -3.19%  '
Set-TextValue $ws 'D48' '2.139.08'
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('E49').Value = '  -7.13%  '
$ws.Range('E50').Value = '  +3.60%  '
Set-TextValue $ws 'D51' '9.16'
$ws.Range('E51').Value = '  +0.83%  '
